$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt1"
$ws.Range("C2").Value = "Fzd3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01948966666666667
$ws.Range("H2").Value = 0.058469
$ws.Range("I2").Value = 0.0709606244933031
$ws.Range("J2").Value = 0.0709606244933031
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2078313333333333
$ws.Range("N2").Value = 0.623494
$ws.Range("O2").Value = 0.08621557350328635
$ws.Range("P2").Value = 0.112461889302165
$ws.Range("Q2").Value = 0.004050563409555555
$ws.Range("R2").Value = 0.036455070686
$ws.Range("S2").Value = 0.006117910936841476
$ws.Range("T2").Value = 0.007980365896578351

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt1"
$ws.Range("C3").Value = "Fzd3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01948966666666667
$ws.Range("H3").Value = 0.058469
$ws.Range("I3").Value = 0.0709606244933031
$ws.Range("J3").Value = 0.0709606244933031
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.498127
$ws.Range("N3").Value = 1.494381
$ws.Range("O3").Value = 0.206640184103479
$ws.Range("P3").Value = 0.2695469573039334
$ws.Range("Q3").Value = 0.009708329187666667
$ws.Range("R3").Value = 0.08737496268900001
$ws.Range("S3").Value = 0.014663316509394
$ws.Range("T3").Value = 0.01912722042055682

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt1"
$ws.Range("C4").Value = "Fzd3"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01948966666666667
$ws.Range("H4").Value = 0.058469
$ws.Range("I4").Value = 0.0709606244933031
$ws.Range("J4").Value = 0.0709606244933031
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.016887
$ws.Range("N4").Value = 0.050661
$ws.Range("O4").Value = 0.007005307459654767
$ws.Range("P4").Value = 0.009137909545139137
$ws.Range("Q4").Value = 0.0003291220009999999
$ws.Range("R4").Value = 0.002962098009
$ws.Range("S4").Value = 0.000497100992104697
$ws.Range("T4").Value = 0.0006484317678863885

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Wnt1"
$ws.Range("C5").Value = "Fzd3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.01948966666666667
$ws.Range("H5").Value = 0.058469
$ws.Range("I5").Value = 0.0709606244933031
$ws.Range("J5").Value = 0.0709606244933031
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.6877555
$ws.Range("N5").Value = 3.375511
$ws.Range("O5").Value = 0.7001389349335798
$ws.Range("P5").Value = 0.6088532438487625
$ws.Range("Q5").Value = 0.03289379210983333
$ws.Range("R5").Value = 0.197362752659
$ws.Range("S5").Value = 0.04968229605496293
$ws.Range("T5").Value = 0.04320460640828155

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt1"
$ws.Range("C6").Value = "Fzd3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.255165
$ws.Range("H6").Value = 0.7654949999999999
$ws.Range("I6").Value = 0.9290393755066968
$ws.Range("J6").Value = 0.9290393755066969
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2078313333333333
$ws.Range("N6").Value = 0.623494
$ws.Range("O6").Value = 0.08621557350328635
$ws.Range("P6").Value = 0.112461889302165
$ws.Range("Q6").Value = 0.05303128217
$ws.Range("R6").Value = 0.47728153953
$ws.Range("S6").Value = 0.08009766256644486
$ws.Range("T6").Value = 0.1044815234055866

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt1"
$ws.Range("C7").Value = "Fzd3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.255165
$ws.Range("H7").Value = 0.7654949999999999
$ws.Range("I7").Value = 0.9290393755066968
$ws.Range("J7").Value = 0.9290393755066969
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.498127
$ws.Range("N7").Value = 1.494381
$ws.Range("O7").Value = 0.206640184103479
$ws.Range("P7").Value = 0.2695469573039334
$ws.Range("Q7").Value = 0.127104575955
$ws.Range("R7").Value = 1.143941183595
$ws.Range("S7").Value = 0.191976867594085
$ws.Range("T7").Value = 0.2504197368833766

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Wnt1"
$ws.Range("C8").Value = "Fzd3"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.255165
$ws.Range("H8").Value = 0.7654949999999999
$ws.Range("I8").Value = 0.9290393755066968
$ws.Range("J8").Value = 0.9290393755066969
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.016887
$ws.Range("N8").Value = 0.050661
$ws.Range("O8").Value = 0.007005307459654767
$ws.Range("P8").Value = 0.009137909545139137
$ws.Range("Q8").Value = 0.004308971354999999
$ws.Range("R8").Value = 0.038780742195
$ws.Range("S8").Value = 0.00650820646755007
$ws.Range("T8").Value = 0.008489477777252749

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Wnt1"
$ws.Range("C9").Value = "Fzd3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.255165
$ws.Range("H9").Value = 0.7654949999999999
$ws.Range("I9").Value = 0.9290393755066968
$ws.Range("J9").Value = 0.9290393755066969
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.6877555
$ws.Range("N9").Value = 3.375511
$ws.Range("O9").Value = 0.7001389349335798
$ws.Range("P9").Value = 0.6088532438487625
$ws.Range("Q9").Value = 0.4306561321575
$ws.Range("R9").Value = 2.583936792945
$ws.Range("S9").Value = 0.6504566388786168
$ws.Range("T9").Value = 0.565648637440481
